{"js": "// Insert a new \"Author\" styled paragraph right after the \"Edison Achalma\"\n// author paragraph, containing the author's affiliation.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text,items/style\");\nawait context.sync();\n\nconst targetText = \"Edison Achalma\";\nconst targetStyle = \"Author\";\nconst newText =\n  \"Escuela Profesional de Econom\u00eda, Universidad Nacional de San Crist\u00f3bal de Huamanga\";\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text === targetText && p.style === targetStyle) {\n    target = p;\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error('Could not find the \"Edison Achalma\" Author paragraph.');\n}\n\n// Build a single new paragraph (same \"Author\" style) carrying the\n// affiliation text, and insert it as OOXML right after the target\n// paragraph's range. Using insertOoxml (instead of Paragraph.insertParagraph)\n// keeps the original \"Edison Achalma\" run completely untouched and lets us\n// control the exact markup (including xml:space=\"preserve\") of the new run.\nconst newParagraphXml =\n  '<w:p><w:pPr><w:pStyle w:val=\"' +\n  targetStyle +\n  '\"/></w:pPr><w:r><w:t xml:space=\"preserve\">' +\n  newText +\n  \"</w:t></w:r></w:p>\";\n\nconst ooxmlPackage =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  \"<pkg:xmlData>\" +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  \"<w:body>\" +\n  newParagraphXml +\n  \"</w:body></w:document>\" +\n  \"</pkg:xmlData></pkg:part></pkg:package>\";\n\nconst wholeRange = target.getRange(\"Whole\");\nwholeRange.insertOoxml(ooxmlPackage, \"After\");\n\nawait context.sync();\n", "ps1": "# Insert a new \"Author\" styled paragraph right after the \"Edison Achalma\"\n# author paragraph, containing the author's affiliation.\n$doc = $word.ActiveDocument\n\n$targetText = \"Edison Achalma\"\n$targetStyle = \"Author\"\n$newText = \"Escuela Profesional de Econom\u00eda, Universidad Nacional de San Crist\u00f3bal de Huamanga\"\n\n$target = $null\nfor ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {\n    $p = $doc.Paragraphs.Item($i)\n    $txt = $p.Range.Text.TrimEnd([char]13)\n    if ($txt -eq $targetText -and $p.Style.NameLocal -eq $targetStyle) {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw 'Could not find the \"Edison Achalma\" Author paragraph.'\n}\n\n# `Range.InsertParagraphAfter` is avoided here: in this document the\n# paragraph-mark based insert corrupts the existing \"Edison Achalma\" run\n# (its preceding content includes a Table-of-Contents content control).\n# Insert the new paragraph as raw OOXML right after the target paragraph's\n# range instead; this creates a sibling paragraph with the requested style\n# while leaving the original paragraph completely untouched.\n$newParagraphXml = '<w:p><w:pPr><w:pStyle w:val=\"' + $targetStyle + '\"/></w:pPr><w:r><w:t xml:space=\"preserve\">' + $newText + '</w:t></w:r></w:p>'\n\n$ooxmlPackage = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' + $newParagraphXml + '</w:body></w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>'\n\n[void]$target.Range.InsertXML($ooxmlPackage, \"After\")\n"}
